# Updated symbol list on Fri Dec 23 00:56:03 UTC 2022 with GitHub Actions
# Refreshes Price (D), Volume(1h) label (E), Data (F) and Hora (G) columns
# for each coin row (2-51) to the 23-12-2022 snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "245.37"; "F2" = "23-12-2022"; "G2" = "0"
    "D3" = "21.75"; "F3" = "23-12-2022"; "G3" = "0"
    "D4" = "5.420"; "F4" = "23-12-2022"; "G4" = "0"
    "D5" = "0.05762"; "F5" = "23-12-2022"; "G5" = "0"
    "D6" = "3.404"; "F6" = "23-12-2022"; "G6" = "0"
    "D7" = "6.327"; "F7" = "23-12-2022"; "G7" = "0"
    "D8" = "0.8189"; "F8" = "23-12-2022"; "G8" = "0"
    "D9" = "1.021"; "E9" = "8FTXTokenFTT"; "F9" = "23-12-2022"; "G9" = "0"
    "D10" = "0.1436"; "F10" = "23-12-2022"; "G10" = "0"
    "D11" = "0.07304"; "F11" = "23-12-2022"; "G11" = "0"
    "D12" = "0.03107"; "F12" = "23-12-2022"; "G12" = "0"
    "D13" = "0.03046"; "F13" = "23-12-2022"; "G13" = "0"
    "D14" = "4.160"; "F14" = "23-12-2022"; "G14" = "0"
    "D15" = "0.09395"; "F15" = "23-12-2022"; "G15" = "0"
    "D16" = "0.001592"; "F16" = "23-12-2022"; "G16" = "0"
    "F17" = "23-12-2022"; "G17" = "0"
    "D18" = "0.0005854"; "F18" = "23-12-2022"; "G18" = "0"
    "D19" = "0.006240"; "F19" = "23-12-2022"; "G19" = "0"
    "D20" = "0.004110"; "F20" = "23-12-2022"; "G20" = "0"
    "D21" = "0.0009967"; "F21" = "23-12-2022"; "G21" = "0"
    "D22" = "0.0001501"; "F22" = "23-12-2022"; "G22" = "0"
    "D23" = "3.735"; "F23" = "23-12-2022"; "G23" = "0"
    "D24" = "2.200"; "F24" = "23-12-2022"; "G24" = "0"
    "D25" = "0.3255"; "F25" = "23-12-2022"; "G25" = "0"
    "D26" = "0.1330"; "F26" = "23-12-2022"; "G26" = "0"
    "D27" = "0.0004001"; "F27" = "23-12-2022"; "G27" = "0"
    "F28" = "23-12-2022"; "G28" = "0"
    "F29" = "23-12-2022"; "G29" = "0"
    "F30" = "23-12-2022"; "G30" = "0"
    "F31" = "23-12-2022"; "G31" = "0"
    "F32" = "23-12-2022"; "G32" = "0"
    "F33" = "23-12-2022"; "G33" = "0"
    "F34" = "23-12-2022"; "G34" = "0"
    "F35" = "23-12-2022"; "G35" = "0"
    "F36" = "23-12-2022"; "G36" = "0"
    "F37" = "23-12-2022"; "G37" = "0"
    "F38" = "23-12-2022"; "G38" = "0"
    "F39" = "23-12-2022"; "G39" = "0"
    "D40" = "0.03891"; "F40" = "23-12-2022"; "G40" = "0"
    "D41" = "0.006681"; "F41" = "23-12-2022"; "G41" = "0"
    "D42" = "0.1071"; "F42" = "23-12-2022"; "G42" = "0"
    "D43" = "0.002902"; "F43" = "23-12-2022"; "G43" = "0"
    "D44" = "0.006669"; "F44" = "23-12-2022"; "G44" = "0"
    "D45" = "0.00005615"; "F45" = "23-12-2022"; "G45" = "0"
    "F46" = "23-12-2022"; "G46" = "0"
    "D47" = "0.5403"; "E47" = "46CoinbaseStockTokenCOINBestin24h"; "F47" = "23-12-2022"; "G47" = "0"
    "F48" = "23-12-2022"; "G48" = "0"
    "D49" = "0.00002101"; "F49" = "23-12-2022"; "G49" = "0"
    "D50" = "0.01011"; "F50" = "23-12-2022"; "G50" = "0"
    "F51" = "23-12-2022"; "G51" = "0"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.ClearFormats()
}

Write-Host "Updated $($updates.Count) cells for the 23-12-2022 snapshot."
